# Nexial mail-showcase.xlsx update
# - [web] adds new command `deselect(locator,text)` to the '#system' sheet's
#   "web" lookup column (column U), inserted alphabetically right before
#   `deselectMulti(locator,array)` (which sits at the former U53).
# - This shifts U53:U116 down to U54:U117 and grows the named range `web`
#   from $U$2:$U$116 to $U$2:$U$117.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("#system")

# Insert a new row above the current U53 ("deselectMulti(locator,array)"),
# pushing it (and everything below it through U116) down by one row.
$ws.Rows.Item(53).Insert()

# Populate the freshly inserted row with the new web-automation command.
$ws.Range("U53").Value = "deselect(locator,text)"

# Grow the "web" defined name so it covers the newly added row.
$wb.Names.Item("web").RefersTo = "='#system'!`$U`$2:`$U`$117"
